$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.71508866666667
$ws.Range("H2").Value = 107.145266
$ws.Range("I2").Value = 0.1390302752364672
$ws.Range("J2").Value = 0.1390302752364672
$ws.Range("M2").Value = 1.234379333333333
$ws.Range("N2").Value = 3.703138
$ws.Range("O2").Value = 0.2458010442471192
$ws.Range("P2").Value = 0.2458010442471192
$ws.Range("Q2").Value = 44.0859673383009
$ws.Range("R2").Value = 396.7737060447081
$ws.Range("S2").Value = 0.03417378683508802
$ws.Range("T2").Value = 0.03417378683508802
$ws.Range("G3").Value = 35.71508866666667
$ws.Range("H3").Value = 107.145266
$ws.Range("I3").Value = 0.1390302752364672
$ws.Range("J3").Value = 0.1390302752364672
$ws.Range("O3").Value = 0.3817123403920895
$ws.Range("P3").Value = 0.3817123403920895
$ws.Range("Q3").Value = 68.46251537578357
$ws.Range("R3").Value = 616.1626383820521
$ws.Range("S3").Value = 0.05306957174586825
$ws.Range("T3").Value = 0.05306957174586826
$ws.Range("G4").Value = 35.71508866666667
$ws.Range("H4").Value = 107.145266
$ws.Range("I4").Value = 0.1390302752364672
$ws.Range("J4").Value = 0.1390302752364672
$ws.Range("M4").Value = 1.411497
$ws.Range("N4").Value = 4.234491
$ws.Range("O4").Value = 0.2810703542927722
$ws.Range("P4").Value = 0.2810703542927722
$ws.Range("Q4").Value = 50.41174050773401
$ws.Range("R4").Value = 453.7056645696061
$ws.Range("S4").Value = 0.03907728871813546
$ws.Range("T4").Value = 0.03907728871813547
$ws.Range("G5").Value = 35.71508866666667
$ws.Range("H5").Value = 107.145266
$ws.Range("I5").Value = 0.1390302752364672
$ws.Range("J5").Value = 0.1390302752364672
$ws.Range("M5").Value = 0.45908
$ws.Range("N5").Value = 1.37724
$ws.Range("O5").Value = 0.09141626106801917
$ws.Range("P5").Value = 0.09141626106801917
$ws.Range("Q5").Value = 16.39608290509334
$ws.Range("R5").Value = 147.56474614584
$ws.Range("S5").Value = 0.01270962793737544
$ws.Range("T5").Value = 0.01270962793737545
$ws.Range("G6").Value = 54.09018966666667
$ws.Range("I6").Value = 0.2105601368412127
$ws.Range("J6").Value = 0.2105601368412127
$ws.Range("M6").Value = 1.234379333333333
$ws.Range("N6").Value = 3.703138
$ws.Range("O6").Value = 0.2458010442471192
$ws.Range("P6").Value = 0.2458010442471192
$ws.Range("Q6").Value = 66.76781226061355
$ws.Range("R6").Value = 600.910310345522
$ws.Range("S6").Value = 0.05175590151238638
$ws.Range("T6").Value = 0.05175590151238638
$ws.Range("G7").Value = 54.09018966666667
$ws.Range("I7").Value = 0.2105601368412127
$ws.Range("J7").Value = 0.2105601368412127
$ws.Range("O7").Value = 0.3817123403920895
$ws.Range("P7").Value = 0.3817123403920895
$ws.Range("R7").Value = 933.1729311008179
$ws.Range("S7").Value = 0.08037340262693793
$ws.Range("T7").Value = 0.08037340262693793
$ws.Range("G8").Value = 54.09018966666667
$ws.Range("I8").Value = 0.2105601368412127
$ws.Range("J8").Value = 0.2105601368412127
$ws.Range("M8").Value = 1.411497
$ws.Range("N8").Value = 4.234491
$ws.Range("O8").Value = 0.2810703542927722
$ws.Range("P8").Value = 0.2810703542927722
$ws.Range("Q8").Value = 76.348140443931
$ws.Range("R8").Value = 687.133263995379
$ws.Range("S8").Value = 0.05918221226189425
$ws.Range("T8").Value = 0.05918221226189425
$ws.Range("G9").Value = 54.09018966666667
$ws.Range("I9").Value = 0.2105601368412127
$ws.Range("J9").Value = 0.2105601368412127
$ws.Range("M9").Value = 0.45908
$ws.Range("N9").Value = 1.37724
$ws.Range("O9").Value = 0.09141626106801917
$ws.Range("P9").Value = 0.09141626106801917
$ws.Range("Q9").Value = 24.83172427217333
$ws.Range("R9").Value = 223.48551844956
$ws.Range("S9").Value = 0.01924862043999414
$ws.Range("T9").Value = 0.01924862043999414
$ws.Range("G10").Value = 101.4529346666666
$ws.Range("H10").Value = 304.358804
$ws.Range("I10").Value = 0.3949319449238378
$ws.Range("J10").Value = 0.3949319449238378
$ws.Range("M10").Value = 1.234379333333333
$ws.Range("N10").Value = 3.703138
$ws.Range("O10").Value = 0.2458010442471192
$ws.Range("P10").Value = 0.2458010442471192
$ws.Range("Q10").Value = 125.2314058585502
$ws.Range("R10").Value = 1127.082652726952
$ws.Range("S10").Value = 0.09707468446882507
$ws.Range("T10").Value = 0.09707468446882508
$ws.Range("G11").Value = 101.4529346666666
$ws.Range("H11").Value = 304.358804
$ws.Range("I11").Value = 0.3949319449238378
$ws.Range("J11").Value = 0.3949319449238378
$ws.Range("O11").Value = 0.3817123403920895
$ws.Range("P11").Value = 0.3817123403920895
$ws.Range("Q11").Value = 194.4758744507209
$ws.Range("R11").Value = 1750.282870056488
$ws.Range("S11").Value = 0.1507503969924779
$ws.Range("T11").Value = 0.1507503969924779
$ws.Range("G12").Value = 101.4529346666666
$ws.Range("H12").Value = 304.358804
$ws.Range("I12").Value = 0.3949319449238378
$ws.Range("J12").Value = 0.3949319449238378
$ws.Range("M12").Value = 1.411497
$ws.Range("N12").Value = 4.234491
$ws.Range("O12").Value = 0.2810703542927722
$ws.Range("P12").Value = 0.2810703542927722
$ws.Range("Q12").Value = 143.200512923196
$ws.Range("R12").Value = 1288.804616308764
$ws.Range("S12").Value = 0.1110036616812767
$ws.Range("T12").Value = 0.1110036616812767
$ws.Range("G13").Value = 101.4529346666666
$ws.Range("H13").Value = 304.358804
$ws.Range("I13").Value = 0.3949319449238378
$ws.Range("J13").Value = 0.3949319449238378
$ws.Range("M13").Value = 0.45908
$ws.Range("N13").Value = 1.37724
$ws.Range("O13").Value = 0.09141626106801917
$ws.Range("P13").Value = 0.09141626106801917
$ws.Range("Q13").Value = 46.57501324677332
$ws.Range("R13").Value = 419.17511922096
$ws.Range("S13").Value = 0.03610320178125812
$ws.Range("T13").Value = 0.03610320178125813
$ws.Range("G14").Value = 65.628919
$ws.Range("H14").Value = 196.886757
$ws.Range("I14").Value = 0.2554776429984823
$ws.Range("J14").Value = 0.2554776429984823
$ws.Range("M14").Value = 1.234379333333333
$ws.Range("N14").Value = 3.703138
$ws.Range("O14").Value = 0.2458010442471192
$ws.Range("P14").Value = 0.2458010442471192
$ws.Range("Q14").Value = 81.01098128260732
$ws.Range("R14").Value = 729.098831543466
$ws.Range("S14").Value = 0.06279667143081966
$ws.Range("T14").Value = 0.06279667143081966
$ws.Range("G15").Value = 65.628919
$ws.Range("H15").Value = 196.886757
$ws.Range("I15").Value = 0.2554776429984823
$ws.Range("J15").Value = 0.2554776429984823
$ws.Range("O15").Value = 0.3817123403920895
$ws.Range("P15").Value = 0.3817123403920895
$ws.Range("Q15").Value = 125.8045561098393
$ws.Range("R15").Value = 1132.241004988554
$ws.Range("S15").Value = 0.0975189690268054
$ws.Range("T15").Value = 0.0975189690268054
$ws.Range("G16").Value = 65.628919
$ws.Range("H16").Value = 196.886757
$ws.Range("I16").Value = 0.2554776429984823
$ws.Range("J16").Value = 0.2554776429984823
$ws.Range("M16").Value = 1.411497
$ws.Range("N16").Value = 4.234491
$ws.Range("O16").Value = 0.2810703542927722
$ws.Range("P16").Value = 0.2810703542927722
$ws.Range("Q16").Value = 92.63502228174299
$ws.Range("R16").Value = 833.715200535687
$ws.Range("S16").Value = 0.0718071916314658
$ws.Range("T16").Value = 0.0718071916314658
$ws.Range("G17").Value = 65.628919
$ws.Range("H17").Value = 196.886757
$ws.Range("I17").Value = 0.2554776429984823
$ws.Range("J17").Value = 0.2554776429984823
$ws.Range("M17").Value = 0.45908
$ws.Range("N17").Value = 1.37724
$ws.Range("O17").Value = 0.09141626106801917
$ws.Range("P17").Value = 0.09141626106801917
$ws.Range("Q17").Value = 30.12892413452
$ws.Range("R17").Value = 271.16031721068
$ws.Range("S17").Value = 0.02335481090939146
$ws.Range("T17").Value = 0.02335481090939146
